$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 195.19048
$ws.Range("I92").Value = 175.45
$ws.Range("J92").Value = 590
$ws.Range("K92").Value = 175.45
$ws.Range("L92").Value = 590
$ws.Range("M92").Value = 1072.55
$ws.Range("N92").Value = -3086

$ws.Range("H107").Value = 1028.6666
$ws.Range("I107").Value = 1042.5714
$ws.Range("K107").Value = 1042.5714
$ws.Range("M107").Value = 877.4286

$ws.Range("H111").Value = 3648.3572
$ws.Range("I111").Value = 2837.7
$ws.Range("K111").Value = 8513.099999999999
$ws.Range("M111").Value = -5446.099999999999

$ws.Range("H137").Value = 2396.1
$ws.Range("I137").Value = 1416.3636
$ws.Range("J137").Value = 2963.3157
$ws.Range("K137").Value = 4249.0908
$ws.Range("L137").Value = 8889.947100000001
$ws.Range("M137").Value = -1699.0908
$ws.Range("N137").Value = -13989.9471

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2852.9524
$ws.Range("I2").Value = 1800.8462
$ws.Range("J2").Value = 4562.625
$ws.Range("K2").Value = 1800.8462
$ws.Range("L2").Value = 4562.625
$ws.Range("M2").Value = -1687.8462
$ws.Range("N2").Value = -4788.625

$ws.Range("H4").Value = 103.333336
$ws.Range("I4").Value = 95
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 95
$ws.Range("L4").Value = 120
$ws.Range("M4").Value = 21
$ws.Range("N4").Value = -352

$ws.Range("H5").Value = 110.44444
$ws.Range("I5").Value = 107.333336
$ws.Range("J5").Value = 116.666664
$ws.Range("K5").Value = 107.333336
$ws.Range("L5").Value = 116.666664
$ws.Range("M5").Value = 4.666663999999997
$ws.Range("N5").Value = -340.666664

$ws.Range("H32").Value = 26518.363
$ws.Range("I32").Value = 9727.763999999999
$ws.Range("J32").Value = 102076.06
$ws.Range("K32").Value = 9727.763999999999
$ws.Range("L32").Value = 102076.06
$ws.Range("M32").Value = -9440.763999999999
$ws.Range("N32").Value = -102650.06

$ws.Range("H61").Value = 2391.4167
$ws.Range("I61").Value = 2373.158
$ws.Range("J61").Value = 2460.8
$ws.Range("K61").Value = 2373.158
$ws.Range("L61").Value = 2460.8
$ws.Range("M61").Value = -2161.158
$ws.Range("N61").Value = -2884.8

$ws.Range("H74").Value = 4636.6
$ws.Range("I74").Value = 4257.9375
$ws.Range("J74").Value = 5309.778
$ws.Range("K74").Value = 4257.9375
$ws.Range("L74").Value = 5309.778
$ws.Range("M74").Value = -3383.9375
$ws.Range("N74").Value = -7057.778

$ws.Range("H77").Value = 4636.6
$ws.Range("I77").Value = 4257.9375
$ws.Range("J77").Value = 5309.778
$ws.Range("K77").Value = 21289.6875
$ws.Range("L77").Value = 26548.89
$ws.Range("M77").Value = -16921.6875
$ws.Range("N77").Value = -35284.89

$ws.Range("H110").Value = 1317.4572
$ws.Range("I110").Value = 1620.7391
$ws.Range("K110").Value = 1620.7391
$ws.Range("M110").Value = 424.2609

$ws.Range("H116").Value = 2852.9524
$ws.Range("I116").Value = 1800.8462
$ws.Range("J116").Value = 4562.625
$ws.Range("K116").Value = 1800.8462
$ws.Range("L116").Value = 4562.625
$ws.Range("M116").Value = 493.1538
$ws.Range("N116").Value = -9150.625

$ws.Range("H132").Value = 3138.1428
$ws.Range("I132").Value = 1954.7307
$ws.Range("J132").Value = 6556.8887
$ws.Range("K132").Value = 5864.1921
$ws.Range("L132").Value = 19670.6661
$ws.Range("M132").Value = -3334.1921
$ws.Range("N132").Value = -24730.6661

$ws.Range("H136").Value = 2391.4167
$ws.Range("I136").Value = 2373.158
$ws.Range("J136").Value = 2460.8
$ws.Range("K136").Value = 7119.474
$ws.Range("L136").Value = 7382.400000000001
$ws.Range("M136").Value = -4569.474
$ws.Range("N136").Value = -12482.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2852.9524
$ws.Range("I3").Value = 1800.8462
$ws.Range("J3").Value = 4562.625
$ws.Range("K3").Value = 1800.8462
$ws.Range("L3").Value = 4562.625
$ws.Range("M3").Value = -1686.8462
$ws.Range("N3").Value = -4790.625

$ws.Range("H4").Value = 110.44444
$ws.Range("I4").Value = 107.333336
$ws.Range("J4").Value = 116.666664
$ws.Range("K4").Value = 107.333336
$ws.Range("L4").Value = 116.666664
$ws.Range("M4").Value = 7.666663999999997
$ws.Range("N4").Value = -346.666664

$ws.Range("H94").Value = 637.30554
$ws.Range("I94").Value = 559.6818
$ws.Range("J94").Value = 759.2857
$ws.Range("K94").Value = 559.6818
$ws.Range("L94").Value = 759.2857
$ws.Range("M94").Value = -108.6818
$ws.Range("N94").Value = -1661.2857

$ws.Range("H134").Value = 2342.5789
$ws.Range("I134").Value = 1835
$ws.Range("J134").Value = 3040.5
$ws.Range("K134").Value = 5505
$ws.Range("L134").Value = 9121.5
$ws.Range("M134").Value = -2970
$ws.Range("N134").Value = -14191.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 72.2
$ws.Range("I7").Value = 41.846153
$ws.Range("J7").Value = 128.57143
$ws.Range("K7").Value = 41.846153
$ws.Range("L7").Value = 128.57143
$ws.Range("M7").Value = 71.153847
$ws.Range("N7").Value = -354.57143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3150.4
$ws.Range("I134").Value = 2727.0588
$ws.Range("K134").Value = 8181.176399999999
$ws.Range("M134").Value = -3111.176399999999

$ws.Range("H139").Value = 14669.8
$ws.Range("I139").Value = 1576.2963
$ws.Range("J139").Value = 21340.076
$ws.Range("K139").Value = 4728.8889
$ws.Range("L139").Value = 64020.228
$ws.Range("M139").Value = 411.1111000000001
$ws.Range("N139").Value = -74300.228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2173.3635
$ws.Range("I132").Value = 2145.1052
$ws.Range("J132").Value = 2352.3333
$ws.Range("K132").Value = 6435.3156
$ws.Range("L132").Value = 7056.999899999999
$ws.Range("M132").Value = -3905.3156
$ws.Range("N132").Value = -12116.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1500.3334
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1750.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1750.5
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2340.5

$ws.Range("H27").Value = 1500.3334
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1750.5
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1750.5
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1964.5

$ws.Range("H46").Value = 2333.3333
$ws.Range("I46").Value = 1633.3334
$ws.Range("J46").Value = 3033.3333
$ws.Range("K46").Value = 1633.3334
$ws.Range("L46").Value = 3033.3333
$ws.Range("M46").Value = -1445.3334
$ws.Range("N46").Value = -3409.3333

$ws.Range("H55").Value = 971.1539
$ws.Range("I55").Value = 255.16667
$ws.Range("J55").Value = 1584.8572
$ws.Range("K55").Value = 255.16667
$ws.Range("L55").Value = 1584.8572
$ws.Range("M55").Value = -82.16667000000001
$ws.Range("N55").Value = -1930.8572

$ws.Range("H68").Value = 3135.2
$ws.Range("I68").Value = 3087.75
$ws.Range("J68").Value = 3325
$ws.Range("K68").Value = 3087.75
$ws.Range("L68").Value = 3325
$ws.Range("M68").Value = -2338.75
$ws.Range("N68").Value = -4823

$ws.Range("H71").Value = 3135.2
$ws.Range("I71").Value = 3087.75
$ws.Range("J71").Value = 3325
$ws.Range("K71").Value = 15438.75
$ws.Range("L71").Value = 16625
$ws.Range("M71").Value = -11694.75
$ws.Range("N71").Value = -24113

$ws.Range("H132").Value = 2658.0344
$ws.Range("I132").Value = 2453.5908
$ws.Range("J132").Value = 3300.5715
$ws.Range("K132").Value = 7360.7724
$ws.Range("L132").Value = 9901.7145
$ws.Range("M132").Value = -4830.7724
$ws.Range("N132").Value = -14961.7145

$ws.Range("H138").Value = 56856
$ws.Range("J138").Value = 56856
$ws.Range("L138").Value = 56856
$ws.Range("N138").Value = -67136

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8152.875
$ws.Range("I81").Value = 14308.25
$ws.Range("J81").Value = 1997.5
$ws.Range("K81").Value = 28616.5
$ws.Range("L81").Value = 3995
$ws.Range("M81").Value = -27555.5
$ws.Range("N81").Value = -6117

$ws.Range("H84").Value = 8152.875
$ws.Range("I84").Value = 14308.25
$ws.Range("J84").Value = 1997.5
$ws.Range("K84").Value = 143082.5
$ws.Range("L84").Value = 19975
$ws.Range("M84").Value = -137778.5
$ws.Range("N84").Value = -30583

$ws.Range("H132").Value = 3020.4167
$ws.Range("J132").Value = 1818.8182
$ws.Range("L132").Value = 5456.4546
$ws.Range("N132").Value = -10516.4546
